$wb = $excel.ActiveWorkbook

# Rename sheets (task-order IDs refreshed)
$wb.Worksheets.Item(1).Name = "GNG_TO-16511687696314583"
$wb.Worksheets.Item(2).Name = "NB_TO-16511687720711088"
$wb.Worksheets.Item(3).Name = "RS_TO-16511687720730774"
$wb.Worksheets.Item(4).Name = "TOL_TO-16511687721201072"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16511687721990879"

# Sheet 1 (GNG) - update stim filenames
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16511687695873928.csv"
$ws1.Range("B3").Value = "GNG_stims-1651168769613583.csv"
$ws1.Range("B4").Value = "go_stims-16511687696155505.csv"
$ws1.Range("B5").Value = "GNG_stims-16511687696306753.csv"

# Sheet 2 (NB) - update stim filenames
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_5-16511687701420736.csv"
$ws2.Range("B3").Value = "TB-16511687704069247.csv"
$ws2.Range("B4").Value = "OB-1651168770165066.csv"
$ws2.Range("B5").Value = "TB-16511687720608735.csv"
$ws2.Range("B6").Value = "ZB-match_1-16511687697035475.csv"
$ws2.Range("B7").Value = "OB-16511687702792926.csv"
$ws2.Range("B8").Value = "TB-1651168771880595.csv"
$ws2.Range("B9").Value = "ZB-match_1-16511687699527915.csv"
$ws2.Range("B10").Value = "OB-1651168770231404.csv"

# Sheet 4 (TOL) - update stim filenames
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16511687720869234.csv"
$ws4.Range("B3").Value = "ZM_stims-16511687720750763.csv"
$ws4.Range("B4").Value = "MM_stims-16511687721025147.csv"
$ws4.Range("B5").Value = "ZM_stims-16511687720869234.csv"
$ws4.Range("B6").Value = "MM_stims-16511687721191254.csv"
$ws4.Range("B7").Value = "ZM_stims-16511687721025147.csv"

# Sheet 5 (vSAT) - update stim filenames
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-1651168772124149.csv"
$ws5.Range("B3").Value = "vSAT_stims-16511687721674101.csv"
$ws5.Range("B4").Value = "vSAT_stims-1651168772183112.csv"
$ws5.Range("B5").Value = "SAT_stims-16511687721504402.csv"
